# The converter that generated this sheet now writes the (previously
# manually-pasted) "include file" values itself, so the placeholder
# columns G:J on the "Main" sheet no longer need their hard-coded
# header labels / sample numbers - clear them out.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Remove the placeholder header labels in G1:J1 (style is kept, value/text is not)
$ws.Range("G1:J1").ClearContents()

# Remove the placeholder sample data in G2:J5 entirely
$ws.Range("G2:J5").ClearContents()

# Leave the selection where the user last clicked
$ws.Range("J14").Select()
